# Update crypto price/volume data per Thu Aug 29 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.843.63'
$ws.Range("E2").Value = '  +1.92%  '

$ws.Range("D3").Value = '2.546.75'
$ws.Range("E3").Value = '  +2.50%  '

$ws.Range("E4").Value = '  +0.28%  '

$ws.Range("D5").Value = "'537.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.26%  '

$ws.Range("D6").Value = "'142.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.06%  '

$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.20%  '

$ws.Range("D8").Value = "'0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("D9").Value = '2.560.03'
$ws.Range("E9").Value = '  +2.08%  '

$ws.Range("D10").Value = "'0.100"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.40%  '

$ws.Range("D11").Value = "'0.160"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.83%  '

$ws.Range("D12").Value = "'5.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.57%  '

$ws.Range("D13").Value = "'0.359"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.22%  '

$ws.Range("D14").Value = '2.996.97'
$ws.Range("E14").Value = '  +2.61%  '

$ws.Range("D15").Value = "'23.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.40%  '

$ws.Range("D16").Value = '59.882.01'
$ws.Range("E16").Value = '  +2.20%  '

$ws.Range("E17").Value = '  +3.84%  '

$ws.Range("D18").Value = '2.543.85'
$ws.Range("E18").Value = '  +1.89%  '

$ws.Range("D19").Value = "'11.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.09%  '

$ws.Range("D20").Value = "'4.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.27%  '

$ws.Range("D21").Value = "'325.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.73%  '

$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.03%  '

$ws.Range("D23").Value = "'5.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.31%  '

$ws.Range("D24").Value = "'63.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.36%  '

$ws.Range("D25").Value = "'0.428"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.74%  '

$ws.Range("D26").Value = "'0.166"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.65%  '

$ws.Range("D27").Value = "'0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.32%  '

$ws.Range("D28").Value = "'7.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.46%  '

$ws.Range("D29").Value = "'6.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.15%  '

$ws.Range("D30").Value = '0.0₃0785'
$ws.Range("E30").Value = '  +3.65%  '

$ws.Range("D31").Value = "'1.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.78%  '

$ws.Range("D32").Value = "'165.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.91%  '

$ws.Range("D33").Value = "'1.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.22%  '

$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.27%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = "'1.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.72%  '

$ws.Range("D36").Value = "'18.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.25%  '

$ws.Range("D37").Value = "'4.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.53%  '

$ws.Range("D38").Value = "'1.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.74%  '

$ws.Range("D39").Value = "'37.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.99%  '

$ws.Range("D40").Value = "'5.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.98%  '

$ws.Range("D41").Value = "'296.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.59%  '

$ws.Range("D42").Value = "'3.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.60%  '

$ws.Range("D43").Value = "'0.823"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.05%  '

$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = "'0.996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.14%  '

$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = "'0.607"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.50%  '

$ws.Range("D46").Value = "'10.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.93%  '

$ws.Range("D47").Value = "'126.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.89%  '

$ws.Range("D48").Value = "'0.0933"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.56%  '

$ws.Range("D49").Value = "'18.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.97%  '

$ws.Range("D50").Value = "'0.0516"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.09%  '

$ws.Range("D51").Value = "'0.0227"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.61%  '
